$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Match Sheet1's view state in the final file (selection moved to H12
# once Sheet2 becomes the active/selected tab).
$sheet1.Range("H12").Select()

# Insert the new worksheet right after Sheet1 (so it becomes sheetId 2 /
# rId2, sits second in the tab strip, and becomes the active sheet).
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "Sheet2"

# --- Row 1: section titles -------------------------------------------------
$ws2.Range("A1").Value = "Tabel Data Training"
$ws2.Range("F1").Value = "Tabel Data Uji"

# --- Row 2: column headers --------------------------------------------------
$ws2.Range("A2").Value = "No"
$ws2.Range("B2").Value = "Kecerahan"
$ws2.Range("C2").Value = "Kejenuhan"
$ws2.Range("D2").Value = "Kelas"
$ws2.Range("F2").Value = "K"
$ws2.Range("G2").Value = "Kecerahan"
$ws2.Range("H2").Value = "Kejenuhan"
$ws2.Range("I2").Value = "Kelas"

# --- Rows 3-9: training data (copied verbatim from Sheet1 A3:D9) -----------
$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = 40
$ws2.Range("C3").Value = 20
$ws2.Range("D3").Value = "Merah"

$ws2.Range("A4").Value = 2
$ws2.Range("B4").Value = 50
$ws2.Range("C4").Value = 50
$ws2.Range("D4").Value = "Biru"

$ws2.Range("A5").Value = 3
$ws2.Range("B5").Value = 60
$ws2.Range("C5").Value = 90
$ws2.Range("D5").Value = "Biru"

$ws2.Range("A6").Value = 4
$ws2.Range("B6").Value = 10
$ws2.Range("C6").Value = 25
$ws2.Range("D6").Value = "Merah"

$ws2.Range("A7").Value = 5
$ws2.Range("B7").Value = 70
$ws2.Range("C7").Value = 70
$ws2.Range("D7").Value = "Biru"

$ws2.Range("A8").Value = 6
$ws2.Range("B8").Value = 60
$ws2.Range("C8").Value = 10
$ws2.Range("D8").Value = "Merah"

$ws2.Range("A9").Value = 7
$ws2.Range("B9").Value = 25
$ws2.Range("C9").Value = 80
$ws2.Range("D9").Value = "Biru"

# --- Row 3, cols F:H: the new "Tabel Data Uji" test point -------------------
$ws2.Range("F3").Value = 5
$ws2.Range("G3").Value = 20
$ws2.Range("H3").Value = 35

# Leave the selection on Sheet2 where Excel's "Move or Copy"/entry flow
# would naturally land it.
$ws2.Range("I3").Select()
